# Auto-generated Excel COM-interop script
# Applies numeric cell updates across 8 worksheets per the target diff.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 1305.375
$ws.Range("I12").Value = 410
$ws.Range("J12").Value = 2200.75
$ws.Range("K12").Value = 410
$ws.Range("L12").Value = 2200.75
$ws.Range("M12").Value = -240
$ws.Range("N12").Value = -2540.75
# Row 29
$ws.Range("H29").Value = 4150
$ws.Range("I29").Value = 300
$ws.Range("K29").Value = 900
$ws.Range("M29").Value = -619
# Row 107
$ws.Range("H107").Value = 1095.0385
$ws.Range("I107").Value = 1049.6364
$ws.Range("J107").Value = 1344.75
$ws.Range("K107").Value = 1049.6364
$ws.Range("L107").Value = 1344.75
$ws.Range("M107").Value = 870.3635999999999
$ws.Range("N107").Value = -5184.75
# Row 111
$ws.Range("H111").Value = 2749
$ws.Range("J111").Value = 2718.6
$ws.Range("L111").Value = 8155.799999999999
$ws.Range("N111").Value = -14289.8

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2768.8906
$ws.Range("I32").Value = 2495.1836
$ws.Range("J32").Value = 3663
$ws.Range("K32").Value = 2495.1836
$ws.Range("L32").Value = 3663
$ws.Range("M32").Value = -2208.1836
$ws.Range("N32").Value = -4237
# Row 74
$ws.Range("H74").Value = 1740.871
$ws.Range("I74").Value = 1546.7826
$ws.Range("K74").Value = 1546.7826
$ws.Range("M74").Value = -672.7826
# Row 77
$ws.Range("H77").Value = 1740.871
$ws.Range("I77").Value = 1546.7826
$ws.Range("K77").Value = 7733.913
$ws.Range("M77").Value = -3365.913
# Row 110
$ws.Range("H110").Value = 4049.0605
$ws.Range("I110").Value = 3800.8572
$ws.Range("J110").Value = 5439
$ws.Range("K110").Value = 3800.8572
$ws.Range("L110").Value = 5439
$ws.Range("M110").Value = -1755.8572
$ws.Range("N110").Value = -9529
# Row 122
$ws.Range("H122").Value = 1417.875
$ws.Range("I122").Value = 1308.8823
$ws.Range("K122").Value = 3926.6469
$ws.Range("M122").Value = -1476.6469
# Row 132
$ws.Range("H132").Value = 2260.8096
$ws.Range("I132").Value = 1971.9
$ws.Range("J132").Value = 2983.0833
$ws.Range("K132").Value = 5915.700000000001
$ws.Range("L132").Value = 8949.249899999999
$ws.Range("M132").Value = -3385.700000000001
$ws.Range("N132").Value = -14009.2499

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2196.077
$ws.Range("I86").Value = 2850.5715
$ws.Range("K86").Value = 2850.5715
$ws.Range("M86").Value = -1727.5715
# Row 89
$ws.Range("H89").Value = 2196.077
$ws.Range("I89").Value = 2850.5715
$ws.Range("K89").Value = 14252.8575
$ws.Range("M89").Value = -8636.8575
# Row 94
$ws.Range("H94").Value = 7937275.5
$ws.Range("I94").Value = 9804461
$ws.Range("J94").Value = 1737
$ws.Range("K94").Value = 9804461
$ws.Range("L94").Value = 1737
$ws.Range("M94").Value = -9804010
$ws.Range("N94").Value = -2639
# Row 138
$ws.Range("H138").Value = 80823.42999999999
$ws.Range("J138").Value = 80823.42999999999
$ws.Range("L138").Value = 80823.42999999999
$ws.Range("N138").Value = -91103.42999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 799.7308
$ws.Range("I22").Value = 409.18182
$ws.Range("J22").Value = 1086.1333
$ws.Range("K22").Value = 409.18182
$ws.Range("L22").Value = 1086.1333
$ws.Range("M22").Value = -59.18182000000002
$ws.Range("N22").Value = -1786.1333
# Row 41
$ws.Range("H41").Value = 18129.4
$ws.Range("I41").Value = 4549.1665
$ws.Range("J41").Value = 38499.75
$ws.Range("K41").Value = 4549.1665
$ws.Range("L41").Value = 38499.75
$ws.Range("M41").Value = -4121.1665
$ws.Range("N41").Value = -39355.75
# Row 58
$ws.Range("H58").Value = 12099.4
$ws.Range("I58").Value = 2749.5
$ws.Range("J58").Value = 26124.25
$ws.Range("K58").Value = 2749.5
$ws.Range("L58").Value = 26124.25
$ws.Range("M58").Value = -2546.5
$ws.Range("N58").Value = -26530.25
# Row 60
$ws.Range("H60").Value = 44217.445
$ws.Range("J60").Value = 44217.445
$ws.Range("L60").Value = 44217.445
$ws.Range("N60").Value = -45239.445
# Row 132
$ws.Range("H132").Value = 6372.0586
$ws.Range("I132").Value = 5884.1816
$ws.Range("K132").Value = 17652.5448
$ws.Range("M132").Value = -15122.5448
# Row 134
$ws.Range("H134").Value = 3744.6667
$ws.Range("I134").Value = 2993.6572
$ws.Range("J134").Value = 7499.7144
$ws.Range("K134").Value = 8980.971600000001
$ws.Range("L134").Value = 22499.1432
$ws.Range("M134").Value = -6445.971600000001
$ws.Range("N134").Value = -27569.1432
# Row 135
$ws.Range("H135").Value = 99900
$ws.Range("J135").Value = 99900
$ws.Range("L135").Value = 99900
$ws.Range("N135").Value = -110040
# Row 136
$ws.Range("H136").Value = 12099.4
$ws.Range("I136").Value = 2749.5
$ws.Range("J136").Value = 26124.25
$ws.Range("K136").Value = 8248.5
$ws.Range("L136").Value = 78372.75
$ws.Range("M136").Value = -5698.5
$ws.Range("N136").Value = -83472.75
# Row 138
$ws.Range("H138").Value = 79966.336
$ws.Range("J138").Value = 79966.336
$ws.Range("L138").Value = 79966.336
$ws.Range("N138").Value = -90246.336
# Row 141
$ws.Range("H141").Value = 29393.092
$ws.Range("J141").Value = 30032.5
$ws.Range("L141").Value = 30032.5
$ws.Range("N141").Value = -40392.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 564.3333
$ws.Range("I7").Value = 936
$ws.Range("J7").Value = 99.75
$ws.Range("K7").Value = 2808
$ws.Range("L7").Value = 299.25
$ws.Range("M7").Value = -2696
$ws.Range("N7").Value = -523.25
# Row 42
$ws.Range("H42").Value = 16727.75
$ws.Range("J42").Value = 16727.75
$ws.Range("L42").Value = 50183.25
$ws.Range("N42").Value = -51251.25
# Row 94
$ws.Range("H94").Value = 13970.125
$ws.Range("I94").Value = 10018
$ws.Range("J94").Value = 16341.4
$ws.Range("K94").Value = 30054
$ws.Range("L94").Value = 49024.2
$ws.Range("M94").Value = -29378
$ws.Range("N94").Value = -50376.2
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 109
$ws.Range("H109").Value = 7637.92
$ws.Range("I109").Value = 2529.8
$ws.Range("J109").Value = 11043.333
$ws.Range("K109").Value = 7589.400000000001
$ws.Range("L109").Value = 33129.999
$ws.Range("M109").Value = -6549.400000000001
$ws.Range("N109").Value = -35209.999

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4718.759
$ws.Range("I70").Value = 4365.7896
$ws.Range("K70").Value = 4365.7896
$ws.Range("M70").Value = -4095.7896
# Row 73
$ws.Range("H73").Value = 4718.759
$ws.Range("I73").Value = 4365.7896
$ws.Range("K73").Value = 4365.7896
$ws.Range("M73").Value = -3429.7896
# Row 97
$ws.Range("H97").Value = 1711.7
$ws.Range("I97").Value = 1475.85
$ws.Range("J97").Value = 2183.4
$ws.Range("K97").Value = 1475.85
$ws.Range("L97").Value = 2183.4
$ws.Range("M97").Value = -979.8499999999999
$ws.Range("N97").Value = -3175.4
# Row 122
$ws.Range("H122").Value = 129248.82
$ws.Range("I122").Value = 203728.6
$ws.Range("J122").Value = 22849.143
$ws.Range("K122").Value = 611185.8
$ws.Range("L122").Value = 68547.429
$ws.Range("M122").Value = -608735.8
$ws.Range("N122").Value = -73447.429
# Row 135
$ws.Range("H135").Value = 95310.44500000001
$ws.Range("J135").Value = 95310.44500000001
$ws.Range("L135").Value = 95310.44500000001
$ws.Range("N135").Value = -105450.445

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 1398.2812
$ws.Range("I55").Value = 305.5
$ws.Range("J55").Value = 2803.2856
$ws.Range("K55").Value = 305.5
$ws.Range("L55").Value = 2803.2856
$ws.Range("M55").Value = -132.5
$ws.Range("N55").Value = -3149.2856
# Row 132
$ws.Range("H132").Value = 10588.75
$ws.Range("I132").Value = 8486.764999999999
$ws.Range("K132").Value = 25460.295
$ws.Range("M132").Value = -22930.295
# Row 136
$ws.Range("H136").Value = 4124.4546
$ws.Range("I136").Value = 3229.4443
$ws.Range("J136").Value = 8152
$ws.Range("K136").Value = 9688.332900000001
$ws.Range("L136").Value = 24456
$ws.Range("M136").Value = -7138.332900000001
$ws.Range("N136").Value = -29556

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4466304
$ws.Range("I81").Value = 5495990
$ws.Range("K81").Value = 10991980
$ws.Range("M81").Value = -10990919
# Row 84
$ws.Range("H84").Value = 4466304
$ws.Range("I84").Value = 5495990
$ws.Range("K84").Value = 54959900
$ws.Range("M84").Value = -54954596
# Row 132
$ws.Range("H132").Value = 1966.6957
$ws.Range("I132").Value = 1816.3429
$ws.Range("K132").Value = 5449.028700000001
$ws.Range("M132").Value = -2919.028700000001
# Row 139
$ws.Range("H139").Value = 68964.914
$ws.Range("J139").Value = 68964.914
$ws.Range("L139").Value = 68964.914
$ws.Range("N139").Value = -79244.914
